$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.872.67"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.90"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.54"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.53"
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.258"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0608"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0914"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.81"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.627.36"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.564"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.20"
$ws.Range("E15").Value = "  +16.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.85"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.879.45"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.18"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.31"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0698"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.79"
$ws.Range("E22").Value = "  +4.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.12"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.33"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.49"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.109"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.58"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0485"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.428.01"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +4.37%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.79"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.53"
$ws.Range("E40").Value = "  +13.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.550"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.822"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0491"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "53.38"
$ws.Range("E45").Value = "  -7.05%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.37"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.773.52"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.95"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("E51").Value = "  +5.17%  "
